$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.230.37'
$ws.Range('E2').Value = '  -1.10%  '

$ws.Range('D3').Value = '3.744.27'
$ws.Range('E3').Value = '  +0.58%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '623.88'
$ws.Range('E5').Value = '  +0.71%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.25'
$ws.Range('E6').Value = '  -0.69%  '

$ws.Range('D7').Value = '3.742.89'
$ws.Range('E7').Value = '  +0.61%  '

$ws.Range('E9').Value = '  -1.12%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.169'
$ws.Range('E10').Value = '  +3.04%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.32'
$ws.Range('E11').Value = '  -4.23%  '

$ws.Range('E12').Value = '  -3.09%  '

$ws.Range('E13').Value = '  +0.78%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000261'
$ws.Range('E14').Value = '  +1.93%  '

$ws.Range('D15').Value = '4.368.70'
$ws.Range('E15').Value = '  +0.63%  '

$ws.Range('D16').Value = '3.745.24'
$ws.Range('E16').Value = '  +0.69%  '

$ws.Range('D17').Value = '70.271.59'
$ws.Range('E17').Value = '  -1.17%  '

$ws.Range('E18').Value = '  -1.52%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.63'
$ws.Range('E19').Value = '  +1.21%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.84'
$ws.Range('E20').Value = '  -0.47%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '506.89'
$ws.Range('E21').Value = '  -2.45%  '

$ws.Range('E22').Value = '  +0.00%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.726'
$ws.Range('E23').Value = '  -2.76%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.59'
$ws.Range('E24').Value = '  +1.76%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.84'
$ws.Range('E25').Value = '  -2.04%  '

$ws.Range('E26').Value = '  +2.42%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '13.18'
$ws.Range('E27').Value = '  -2.89%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000136'
$ws.Range('E28').Value = '  +22.06%  '

$ws.Range('E29').Value = '  -0.01%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.50'
$ws.Range('E30').Value = '  -1.57%  '

$ws.Range('E31').Value = '  +1.93%  '

$ws.Range('E32').Value = '  -2.51%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '31.30'
$ws.Range('E33').Value = '  -2.44%  '

$ws.Range('E34').Value = '  -0.35%  '

$ws.Range('E35').Value = '  +0.03%  '

$ws.Range('E36').Value = '  +2.09%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.19'
$ws.Range('E37').Value = '  +0.96%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.139'
$ws.Range('E38').Value = '  +3.75%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.341'
$ws.Range('E39').Value = '  -1.74%  '

$ws.Range('E40').Value = '  -6.33%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '50.34'
$ws.Range('E41').Value = '  -2.66%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '46.14'
$ws.Range('E42').Value = '  +2.81%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '433.57'
$ws.Range('E43').Value = '  +0.37%  '

$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.91'
$ws.Range('E44').Value = '  +1.60%  '

$ws.Range('B45').Value = 'Cosmos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.75'
$ws.Range('E45').Value = '  -1.12%  '

$ws.Range('D46').Value = '3.015.42'
$ws.Range('E46').Value = '  -4.55%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0366'
$ws.Range('E47').Value = '  -0.36%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '27.57'
$ws.Range('E48').Value = '  -2.18%  '

$ws.Range('E49').Value = '  -0.04%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '137.93'
$ws.Range('E50').Value = '  -2.11%  '

$ws.Range('E51').Value = '  +0.83%  '
